# Add a new enum ("unit.GroupType") with two items (Hero=0, Enemy=1)
# to the __enums__ worksheet, following the layout convention already
# used by the sheet (an "##var" row describing the enum itself, plus
# "*items" rows listing name/value/comment for each enum member).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- ##var row describing the new enum --------------------------------
$ws.Range("B4").Value2 = "unit.GroupType"   # full_name
$ws.Range("C4").Value2 = $false             # flags
$ws.Range("D4").Value2 = $true              # unique

# -- *items rows: name / value / comment ------------------------------
# Written in this order (Enemy row before Hero row) so that the shared
# string table ends up with the same ordering as produced by Excel.
$ws.Range("H5").Value2 = "Enemy"
$ws.Range("J5").Value2 = 1
$ws.Range("K5").Value2 = "敌人"

$ws.Range("H4").Value2 = "Hero"
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = "英雄"

# -- apply the same font styling Excel applies to freshly typed cells --
$newCells = @("B4", "H4", "K4", "H5", "K5")
foreach ($addr in $newCells) {
    $cellFont = $ws.Range($addr).Font
    $cellFont.Name = "等线"
    $cellFont.Size = 11
}

# -- leave the active selection on the last edited cell ----------------
$ws.Range("K5").Select()
